$wb = $excel.ActiveWorkbook

# --- Sheet 1: Diluted Shares ---
# Update the note for the "Oct 2025 Pre-Funded Warrants (fully exercised)" row:
# the warrant-offering note is clarified to mention the share sale terms.
$wsDiluted = $wb.Worksheets.Item("Diluted Shares")
$wsDiluted.Range("D8").Value = "Sold along with ~140M shares for `$0.70 in the Oct 2025 `$138M offering; "

# --- Sheet 2: rNPV Model ---
# FSCD US Peak Market Penetration bumped from 15% to 16%, and the supporting
# note updated to reflect the new symptomatic/active penetration blend (20% -> 22%).
$wsRnpv = $wb.Worksheets.Item("rNPV Model")
$wsRnpv.Range("B30").Value = 0.16
$wsRnpv.Range("C30").Value = "first & only therapy; strong first-mover; blending symptomatic/Active penetration (22%) with Post-op maintenance penetration (10%)"

# --- Sheet 3: Sensitivity ---
# Bear-case note's immediate cash value per share updated from $0.98 to $0.85.
$wsSens = $wb.Worksheets.Item("Sensitivity")
$wsSens.Range("H20").Value = "Zero fibrosis signals, FSCD discontinued, and UC becomes a commercial flop (immediate cash value is still `$0.85)"

# --- View state updates ---
$wsRnpv.Application.ActiveWindow
$wsRnpv.Activate()
$wsRnpv.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 69

$wsSens.Activate()
$wsSens.Range("H21").Select()
